$d = $word.ActiveDocument

# --- Simple text replacements ---
$ok0 = $d.Content.Find.Execute('I compensated through after-college doubt-solving sessions', $true, $false, $false, $false, $false, $true, 1, $false, 'I compensated through after hours doubt-solving sessions', 2)
if (-not $ok0) { Write-Output "WARNING: replacement 0 not found" }
$ok1 = $d.Content.Find.Execute('Having learned some of these concepts beforehand in junior college, I was able to share some of knowledge with my classmates which helped me gain a deeper understanding of these subjects.', $true, $false, $false, $false, $false, $true, 1, $false, 'Having learned some of these concepts beforehand in junior college, helped me gain a deeper understanding of these subjects.', 2)
if (-not $ok1) { Write-Output "WARNING: replacement 1 not found" }
$ok2 = $d.Content.Find.Execute('further honing my skills and cultivating a genuine interest for practical application in real world scenarios.', $true, $false, $false, $false, $false, $true, 1, $false, 'further honing my skills and cultivating a genuine interest to apply these in real world scenario.', 2)
if (-not $ok2) { Write-Output "WARNING: replacement 2 not found" }
$ok3 = $d.Content.Find.Execute('due to the pandemic. Our first significant client, Cyberking Capitals, an investment consultancy company, entrusted us with their project. In my role as the system architect, I was tasked with ensuring', $true, $false, $false, $false, $false, $true, 1, $false, 'due to the pandemic. In my role as the system architect for our first significant client CyberKing Capitals an investment consultancy, I was tasked with ensuring', 2)
if (-not $ok3) { Write-Output "WARNING: replacement 3 not found" }
$ok4 = $d.Content.Find.Execute('we encountered both technical and non-technical challenges.', $true, $false, $false, $false, $false, $true, 1, $false, 'we encountered various challenges.', 2)
if (-not $ok4) { Write-Output "WARNING: replacement 4 not found" }
$ok5 = $d.Content.Find.Execute('my experience being a part of a start-up environment focused mainly on web technology, my current role', $true, $false, $false, $false, $false, $true, 1, $false, 'my experience being a part of a start-up environment, my current role', 2)
if (-not $ok5) { Write-Output "WARNING: replacement 5 not found" }
$ok6 = $d.Content.Find.Execute('the traceabilility of medical records', $true, $false, $false, $false, $false, $true, 1, $false, 'the traceability of medical records', 2)
if (-not $ok6) { Write-Output "WARNING: replacement 6 not found" }

# --- Soccer paragraph restructure (move sentence, remove tab, wording tweaks) ---
$findRng = $d.Content.Duplicate
$foundPara = $findRng.Find.Execute("Despite being very active academically", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundPara) { Write-Output "WARNING: soccer paragraph not found" }
$socPara = $findRng.Paragraphs(1).Range
$socTarget = $d.Range($socPara.Start, $socPara.End - 1)
$socTarget.Text = 'Despite being very active academically, I was equally passionate about sports, especially soccer. Being a very competitive person, the game provided me with an outlet for my competitive spirit, and allowed me to push and exceed my boundaries and fostered personal growth. As an introverted individual, stepping onto the field and playing in front of a crowd and securing victories helped me bolster my confidence, and taught me the transformative power of stepping out of one’s comfort zone. As a part of the team, we achieved many milestones, finishing as runner-up in the AIT Sports Fest, a state-level inter-collegiate tournament, and clinching victory in the Shahu Trophy, another prestigious state-level soccer tournament organised by AISSMS College being some of the notable ones. These experiences reinforced the notion that teamwork, determination, self-improvement and persistence are integral components of success, both on and off the field. These principles have stayed with me and transcended the field, influencing various aspects of my life.'

# --- Move _GoBack bookmark to end of soccer paragraph ---
$findRng2 = $d.Content.Duplicate
$null = $findRng2.Find.Execute("Despite being very active academically", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$socPara2 = $findRng2.Paragraphs(1).Range
$bmPos = $d.Range($socPara2.End - 1, $socPara2.End - 1)
$null = $d.Bookmarks.Add("_GoBack", $bmPos)

Write-Output "All edits applied"
